$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "device" rows (mac-address / document-type test rows) appended after
# the existing data, mirroring the pattern of the preceding 5-row blocks:
# Finger Print Scanner / IRIS Scanner / Web Camera / Document Scanner / Printer.
$newRows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; DspecId = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; DspecId = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";            Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; DspecId = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";      Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; DspecId = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; DspecId = 920 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Mac
    $ws.Cells.Item($row, 4).Value = $r.Serial
    $ws.Cells.Item($row, 6).Value = $r.DspecId
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = $true
    $ws.Cells.Item($row, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($row, 9).Value = "superadmin"
    $ws.Cells.Item($row, 10).Value = "now()"
    $ws.Cells.Item($row, 11).Value = "now()"
}

# Trailing blank formatted rows (style only, no content) left by the author
# after the new data block.
for ($row = 162; $row -le 166; $row++) {
    $ws.Cells.Item($row, 8).HorizontalAlignment = -4131
}

# Scroll the view down to the newly-added rows and select the next empty
# name cell, matching where the author's cursor ended up.
$excel.ActiveWindow.ScrollRow = 154
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E159").Select() | Out-Null
